$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet Sheet1 -> CALCULATION
$ws.Name = "CALCULATION"

# 2) Insert a new "Total Management Fee" summary row above the existing
#    "Total PPN 11% / Total VAT / Total HC" rows (AQ/AR, rows 11-13),
#    pushing those three rows down to 12-14 and adding a brand-new row 14
#    for "Total HC". The PPN value (previously in AR12) is cleared out
#    (to be filled in manually later) while the old AR11 figure becomes
#    the new "Total Management Fee" total.

# Capture the existing values before they get overwritten.
$origAR11 = $ws.Range("AR11").Value()   # 1719196.80376576 (was "Total PPN 11%")
$origAR12 = $ws.Range("AR12").Value()   # 312581.23704832  (was "Total VAT")
$origAR13 = $ws.Range("AR13").Value()   # 32                (was "Total HC")

# Clone the formatting of row 13 (label style 99 / value style 100) down
# onto the brand new row 14 before row 13's own formatting gets changed.
$ws.Range("AQ13:AR13").Copy()
$ws.Range("AQ14:AR14").PasteSpecial(-4122)  # xlPasteFormats

# Row 14 <- old row 13 ("Total HC")
$ws.Range("AQ14").Value = "Total HC"
$ws.Range("AR14").Value = $origAR13

# Row 13 <- old row 12 ("Total VAT"); also matches the accounting number
# format used by the rows above it (style 98 instead of the old style 100).
$ws.Range("AQ13").Value = "Total VAT"
$ws.Range("AR12").Copy()
$ws.Range("AR13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AR13").Value = $origAR12

# Row 12 <- old row 11 ("Total PPN 11%"), value left blank for now.
$ws.Range("AQ12").Value = "Total PPN 11%"
$ws.Range("AR12").Value = ""

# Row 11 <- new "Total Management Fee" row, keeping the figure that used
# to sit there.
$ws.Range("AQ11").Value = "Total Management Fee"
$ws.Range("AR11").Value = $origAR11

# 3) Restore the active selection to AQ18.
$ws.Range("AQ18").Select() | Out-Null
